$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a numeric-looking text value into a cell while preserving
# its original (default) style - Excel COM auto-converts plain numeric
# strings assigned via .Value into real numbers, so we briefly force the
# cell to Text format, assign, then restore the original style.
function Set-TextValue($cell, $text) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") '245.38'

Set-TextValue $ws.Range("D3") '23.80'

Set-TextValue $ws.Range("D4") '5.342'

Set-TextValue $ws.Range("D5") '0.05834'

Set-TextValue $ws.Range("D6") '6.475'

Set-TextValue $ws.Range("D7") '3.361'

Set-TextValue $ws.Range("D9") '0.9217'

Set-TextValue $ws.Range("D11") '0.07351'

Set-TextValue $ws.Range("D12") '0.03070'

Set-TextValue $ws.Range("D13") '0.03080'

Set-TextValue $ws.Range("D14") '0.09371'

Set-TextValue $ws.Range("D15") '3.859'

Set-TextValue $ws.Range("D16") '0.001560'

Set-TextValue $ws.Range("D17") '0.04694'

$ws.Range("B18").Value = 'One'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextValue $ws.Range("D18") '0.0006015'
$ws.Range("E18").Value = '17OneONE'

$ws.Range("B19").Value = 'TigerCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue $ws.Range("D19") '0.005960'
$ws.Range("E19").Value = '18TigerCashTCH'

$ws.Range("B20").Value = 'BitKan'
$ws.Range("C20").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
Set-TextValue $ws.Range("D20") '0.001244'
$ws.Range("E20").Value = '19BitKanKAN'

$ws.Range("B21").Value = 'HotbitToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
Set-TextValue $ws.Range("D21") '0.004687'
$ws.Range("E21").Value = '20HotbitTokenHTB'

$ws.Range("B22").Value = 'NitroEx'
$ws.Range("C22").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
Set-TextValue $ws.Range("D22") '0.00008819'
$ws.Range("E22").Value = '21NitroExNTXBestin24h'

$ws.Range("B23").Value = 'LEO'
$ws.Range("C23").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws.Range("D23") '3.593'
$ws.Range("E23").Value = '22LEOLEO'

$ws.Range("B24").Value = 'BTSEToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue $ws.Range("D24") '2.158'
$ws.Range("E24").Value = '23BTSETokenBTSE'

Set-TextValue $ws.Range("D25") '0.3229'

Set-TextValue $ws.Range("D40") '0.03847'

$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextValue $ws.Range("D41") '0.1066'
$ws.Range("E41").Value = '40BKEXTokenBKK'

$ws.Range("B42").Value = 'CEJI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
Set-TextValue $ws.Range("D42") '0.002726'
$ws.Range("E42").Value = '41CEJICEJI'

$ws.Range("B43").Value = 'KickToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextValue $ws.Range("D43") '0.003101'
$ws.Range("E43").Value = '42KickTokenKICKWorstin24h'

Set-TextValue $ws.Range("D44") '0.008502'

Set-TextValue $ws.Range("D45") '0.00005261'

Set-TextValue $ws.Range("D47") '0.6536'

Set-TextValue $ws.Range("D48") '0.001865'
$ws.Range("E48").Value = '47BOLOBOLO'

Set-TextValue $ws.Range("D50") '0.0002005'

Write-Output "applied"
